$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the remaining worker record to period 2508 (replacing the old 2507 row
# that is being removed below) and reuse it as the single surviving entry.
$ws.Range("E16").Value = "2508"
$ws.Range("E16").HorizontalAlignment = -4108

# Remove the old period row (2508) entirely; rows below shift up by one.
$ws.Rows("17:17").Delete()

# Update summary figures to reflect a single remaining period/entry.
$ws.Range("E11").Value = 56940
$ws.Range("F13").Value = 1
